$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A13").Value = 1009
$ws.Range("B13").Value = "test 9"
$ws.Range("C13").Value = 525
$ws.Range("D13").Value = "PRJ-2"

[void]$ws.Range("A14").Select()
